# Update the "客単価" (column H) values on the "ABC分析_客構成" sheet:
# each numeric value is re-derived as (current H value) / (count_客構成, column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $hCell = $ws.Cells.Item($r, 8)

    $eVal = $eCell.Value()
    $hVal = $hCell.Value()

    if (($eVal -isnot [string]) -and ($hVal -isnot [string]) -and ($eVal -ne $null) -and ($hVal -ne $null) -and ($eVal -ne 0)) {
        $hCell.Value = $hVal / $eVal
    }
}
